{"js": "// Map of old \"dividend\u00f7divisor=\" text to the new replacement text.\n// All old values are unique within the document, so a direct search\n// and replace per pair is safe and unambiguous.\nconst replacements = [\n  [\"504\u00f73=\", \"173\u00f77=\"],\n  [\"551\u00f78=\", \"328\u00f72=\"],\n  [\"285\u00f72=\", \"941\u00f79=\"],\n  [\"756\u00f79=\", \"880\u00f76=\"],\n  [\"794\u00f73=\", \"643\u00f77=\"],\n  [\"826\u00f78=\", \"726\u00f77=\"],\n  [\"254\u00f72=\", \"227\u00f75=\"],\n  [\"398\u00f78=\", \"107\u00f73=\"],\n  [\"222\u00f76=\", \"986\u00f78=\"],\n  [\"620\u00f72=\", \"930\u00f75=\"],\n  [\"499\u00f77=\", \"824\u00f77=\"],\n  [\"892\u00f76=\", \"827\u00f76=\"],\n  [\"277\u00f78=\", \"753\u00f74=\"],\n  [\"944\u00f72=\", \"497\u00f72=\"],\n  [\"383\u00f72=\", \"220\u00f78=\"],\n  [\"519\u00f79=\", \"269\u00f74=\"],\n  [\"227\u00f79=\", \"681\u00f75=\"],\n  [\"379\u00f78=\", \"364\u00f72=\"],\n  [\"468\u00f73=\", \"804\u00f73=\"],\n  [\"120\u00f74=\", \"453\u00f74=\"],\n  [\"530\u00f75=\", \"175\u00f75=\"],\n  [\"749\u00f74=\", \"437\u00f79=\"],\n  [\"612\u00f72=\", \"448\u00f73=\"],\n  [\"738\u00f73=\", \"425\u00f73=\"],\n  [\"979\u00f73=\", \"509\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Map of old \"dividend\u00f7divisor=\" text to the new replacement text.\n# All old values are unique within the document, so a direct\n# Find/Replace per pair is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"504\u00f73=\", \"173\u00f77=\"),\n  @(\"551\u00f78=\", \"328\u00f72=\"),\n  @(\"285\u00f72=\", \"941\u00f79=\"),\n  @(\"756\u00f79=\", \"880\u00f76=\"),\n  @(\"794\u00f73=\", \"643\u00f77=\"),\n  @(\"826\u00f78=\", \"726\u00f77=\"),\n  @(\"254\u00f72=\", \"227\u00f75=\"),\n  @(\"398\u00f78=\", \"107\u00f73=\"),\n  @(\"222\u00f76=\", \"986\u00f78=\"),\n  @(\"620\u00f72=\", \"930\u00f75=\"),\n  @(\"499\u00f77=\", \"824\u00f77=\"),\n  @(\"892\u00f76=\", \"827\u00f76=\"),\n  @(\"277\u00f78=\", \"753\u00f74=\"),\n  @(\"944\u00f72=\", \"497\u00f72=\"),\n  @(\"383\u00f72=\", \"220\u00f78=\"),\n  @(\"519\u00f79=\", \"269\u00f74=\"),\n  @(\"227\u00f79=\", \"681\u00f75=\"),\n  @(\"379\u00f78=\", \"364\u00f72=\"),\n  @(\"468\u00f73=\", \"804\u00f73=\"),\n  @(\"120\u00f74=\", \"453\u00f74=\"),\n  @(\"530\u00f75=\", \"175\u00f75=\"),\n  @(\"749\u00f74=\", \"437\u00f79=\"),\n  @(\"612\u00f72=\", \"448\u00f73=\"),\n  @(\"738\u00f73=\", \"425\u00f73=\"),\n  @(\"979\u00f73=\", \"509\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute(\n    [ref]$find.Text,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]$find.Replacement.Text,\n    [ref]2\n  )\n}\n"}
